$wb = $excel.ActiveWorkbook

# Rename "Hoja2" to "Suspensiones"
$wsSusp = $wb.Worksheets.Item("Hoja2")
$wsSusp.Name = "Suspensiones"

# Make "Suspensiones" the active/selected tab (was "Horas") so the
# workbook's activeTab / per-sheet tabSelected flags move accordingly.
$wsSusp.Activate()
